$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename shared-string header cell values: _old -> _FV2404, _new -> _FV2410
$usedRange = $ws.Range("A1:U1")
foreach ($cell in $usedRange.Cells) {
    $val = $cell.Value()
    if ($val -ne $null) {
        $text = $val.ToString()
        if ($text.EndsWith("_old")) {
            $cell.Value = $text.Replace("_old", "_FV2404")
        } elseif ($text.EndsWith("_new")) {
            $cell.Value = $text.Replace("_new", "_FV2410")
        }
    }
}

# 2. Freeze first row (pane split) and set selection
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3. Convert range into an Excel Table ("ListObject")
$tableRange = $ws.Range("A1:U70")
$listObject = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$listObject.Name = "Table1"
$listObject.TableStyle = ""
